$wb = $excel.ActiveWorkbook

# --- 1. "ISIC to BLS Map" sheet: split "ISIC 20T21" / "Chemicals and pharmaceutical
#        products" (row 11) into two rows: "ISIC 20" / "Chemicals" (row 11) and
#        "ISIC 21" / "Pharmaceuticals" (row 12), both mapped to the same BLS category
#        ("Nondurable goods") as the original combined row.
$wsMap = $wb.Worksheets.Item("ISIC to BLS Map")

# Insert a new row above the old row 12 ("ISIC 22"), pushing everything below down by one.
$wsMap.Rows.Item(12).Insert()

# Newly inserted row 12 gets the "ISIC 21" code now; its text columns are filled in
# below (after row 11's), with the same mapped-category formula/lookup pattern as the
# row above it.
$wsMap.Range("A12").Value = "ISIC 21"
$wsMap.Range("C12").Formula = '=''BLS Table 3''!$A$42'
$wsMap.Range("D12").Formula = '=VLOOKUP($C12,''BLS Table 3''!$A$3:$F$75,COLUMN(''BLS Table 3''!$F$2),FALSE)/100'

# Row 11 keeps the formulas it already had; just relabel it to the "ISIC 20" split.
$wsMap.Range("B11").Value = "Chemicals"
$wsMap.Range("B12").Value = "Pharmaceuticals"
$wsMap.Range("A11").Value = "ISIC 20"

# --- 2. "URPbIC" sheet: same split, but as columns. Column K held "ISIC 20T21"; it
#        becomes "ISIC 20", and a new column L becomes "ISIC 21".
$wsUrp = $wb.Worksheets.Item("URPbIC")

$wsUrp.Columns.Item(11).Insert()

$wsUrp.Range("K1").Value = "ISIC 20"
$wsUrp.Range("K2").Formula = '=VLOOKUP(K$1,''ISIC to BLS Map''!$A$1:$D$38,COLUMN(''ISIC to BLS Map''!$D$1),FALSE)'

$wsUrp.Range("L1").Value = "ISIC 21"
$wsUrp.Range("L2").Formula = '=VLOOKUP(L$1,''ISIC to BLS Map''!$A$1:$D$38,COLUMN(''ISIC to BLS Map''!$D$1),FALSE)'
